$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 700
$ws.Range("G2").Value = 300
$ws.Range("I2").Value = 700
$ws.Range("K2").Value = 850

$ws.Range("E3").Value = 700
$ws.Range("G3").Value = 300
$ws.Range("I3").Value = 700
$ws.Range("K3").Value = 850

$ws.Range("I4").Select()
